$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("I4").Value = 4
$ws.Range("J4").Value = 5
[void]$ws.Range("J6").Select()
